# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (per commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 0
$wsExhibit.Range("F3").Value = 0
$wsExhibit.Range("F4").Value = 1627
$wsExhibit.Range("F5").Value = 17
$wsExhibit.Range("F8").Value = 144
$wsExhibit.Range("F9").Value = 0
$wsExhibit.Range("F10").Value = 488

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 388
$wsAll.Range("F3").Value = 0
$wsAll.Range("F4").Value = 0
$wsAll.Range("F5").Value = 17
$wsAll.Range("F6").Value = 23
$wsAll.Range("F7").Value = 0
$wsAll.Range("F10").Value = 0
